# Updates the cryptos worksheet with refreshed price / 1h-volume figures,
# and swaps the Maker/Dai rows (row 40 becomes Maker, row 41 becomes Dai)
# to reflect the latest ranking pulled by the scraping GitHub Action.
#
# Numeric-looking price strings (e.g. "587.81") are written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing inline-string convention) instead of silently converting them
# to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.764.78"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "3.564.08"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'587.81"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'189.43"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("D7").Value = "3.554.06"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("D8").Value = "'0.622"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +7.94%  "
$ws.Range("D11").Value = "'0.645"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "'54.16"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "'0.0000310"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "'9.43"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "4.124.32"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "70.745.66"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "3.586.17"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "'12.72"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "'18.95"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'561.67"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'17.92"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").Value = "'4.62"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'4.90"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'93.88"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'11.11"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'2.92"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "'9.32"
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").Value = "'32.36"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "'7.07"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "'12.20"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "'3.97"
$ws.Range("E33").Value = "  +28.62%  "
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").Value = "'63.16"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'3.24"
$ws.Range("E36").Value = "  +6.01%  "
$ws.Range("D37").Value = "'528.78"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.406"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "'38.14"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.649.98"
$ws.Range("E40").Value = "  +9.95%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "0.0₃0787"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("D43").Value = "'3.56"
$ws.Range("E43").Value = "  +6.11%  "
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "'0.0458"
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("D46").Value = "'3.46"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'2.93"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'0.138"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").Value = "'9.18"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'1.46"
$ws.Range("E51").Value = "  +9.58%  "
